$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "round-off" work-experience break date values that were
# mistakenly populated with 31940 in these cells, while keeping their
# existing number-format style.
$ws.Range("DS2:DT2").ClearContents()
$ws.Range("DV2:DW2").ClearContents()
$ws.Range("DY2:DZ2").ClearContents()
$ws.Range("EB2:EC2").ClearContents()

# Update the view/selection state of the sheet to reflect where the user
# was working (scrolled further right, selection spanning DS2:EC2).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 119
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("DS2:EC2").Select() | Out-Null
